# Fix bug superhost e nomenclatura

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header F1: is_superhost -> superhost
$ws.Range("F1").Value = "superhost"

# 2) Rows 3 and 4 had their listing details (subtitle/bedrooms/price/
#    rating/superhost flag) swapped by mistake - swap columns B:F back
#    between the two rows (column A - the title - is identical already)
$colsToSwap = @("B", "C", "D", "E", "F")
foreach ($col in $colsToSwap) {
    $cellRow3 = $ws.Range($col + "3")
    $cellRow4 = $ws.Range($col + "4")
    $val3 = $cellRow3.Value2
    $val4 = $cellRow4.Value2

    if ($val4) { $cellRow3.Value = $val4 } else { $cellRow3.ClearContents() }
    if ($val3) { $cellRow4.Value = $val3 } else { $cellRow4.ClearContents() }
}

# 3) Rows 9 and 10 had their title (column A) swapped by mistake
$a9 = $ws.Range("A9").Value2
$a10 = $ws.Range("A10").Value2
$ws.Range("A9").Value = $a10
$ws.Range("A10").Value = $a9

# 4) Clean up the "Superhost" badge text - it used to contain raw HTML
#    markup (a leftover from scraping); replace it with the plain word
#    "Superhost" everywhere it appears in column F
$lastRow = $ws.Cells($ws.Rows.Count, "F").End(-4162).Row  # xlUp
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("F" + $r)
    $val = $cell.Value2
    if ($val -and $val.ToString().Contains("Superhost")) {
        $cell.Value = "Superhost"
    }
}
